# Apply "Add a partial readme" edits to the "Test" sheet of the workbook.
# This fills in previously-empty Win/Lose prediction cells (column H) and a
# few missing numeric stats (columns D and F) for several rows, and moves
# the active cell selection to H18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")
$ws.Activate()

# --- Column H: fill in "Who Win" predictions ------------------------------
# Rows 17 and 21 get a "Lose" prediction using the existing (unshaded) style.
$ws.Range("H17").Value = "Lose"
$ws.Range("H21").Value = "Lose"

# Rows 18 and 20 get a "Win" prediction, matching the shaded style used for
# the other "Win" cells in the sheet (e.g. H2), so copy that formatting over
# before setting the value.
$ws.Range("H2").Copy()
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("H18").Value = "Win"

$ws.Range("H2").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H20").Value = "Win"

$excel.CutCopyMode = 0

# --- Columns D and F: fill in missing numeric stats -----------------------
$ws.Range("D27").Value = 4
$ws.Range("F27").Value = 3

$ws.Range("D28").Value = 7
$ws.Range("F28").Value = 2

$ws.Range("D30").Value = 6
$ws.Range("F30").Value = 4

$ws.Range("D31").Value = 7
$ws.Range("F31").Value = 1

# --- Move the active selection to H18 -------------------------------------
$ws.Range("H18").Select()
